# Apply row data swaps as described by the diff:
#  - Row 3 and Row 7 swap their full record content (columns A,B,D,E,F,G,H,M)
#  - Row 4 and Row 5 swap their A, Q, R values (rest of the row content is identical)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap row 3 <-> row 7 (columns A, B, D, E, F, G, H) ---
$cols37 = @("A", "B", "D", "E", "F", "G", "H")
foreach ($col in $cols37) {
    $addr3 = "$col" + "3"
    $addr7 = "$col" + "7"
    $tmp = $ws.Range($addr3).Value2
    $ws.Range($addr3).Value2 = $ws.Range($addr7).Value2
    $ws.Range($addr7).Value2 = $tmp
}

# Column M: row3 had "färska gnagspår", row7 had nothing -> after swap row3 empty, row7 has it
$ws.Range("M7").Value2 = $ws.Range("M3").Value2
$ws.Range("M3").Value2 = ""

# --- Swap row 4 <-> row 5 (columns A, Q, R) ---
$cols45 = @("A", "Q", "R")
foreach ($col in $cols45) {
    $addr4 = "$col" + "4"
    $addr5 = "$col" + "5"
    $tmp = $ws.Range($addr4).Value2
    $ws.Range($addr4).Value2 = $ws.Range($addr5).Value2
    $ws.Range($addr5).Value2 = $tmp
}
